# Country "Palestine" is removed from the alphabetically sorted list of
# origin countries on the active sheet (Tabelle1). It currently lives in
# row 43 (between "Pakistan" and "Romania"); deleting the whole row shifts
# every following row up by one and removes the now-unused shared string.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(43).Delete()

# Reproduce the author's final view state: scrolled down with B39 selected.
$ws.Range("B39").Select()
